$wb = $excel.ActiveWorkbook

# --- Update "Last Updated" timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 10:37 AM"

# --- Update the "1 Year" (column F) figures on the Industry Analysis sheet ---
$ws = $wb.Worksheets.Item("Industry Analysis")
$ws.Range("F2").Value = 18.476
$ws.Range("F3").Value = -7.7404
$ws.Range("F4").Value = 30.7972
$ws.Range("F5").Value = -50.2266
$ws.Range("F6").Value = 61.9649
$ws.Range("F7").Value = -9.1713
$ws.Range("F8").Value = -3.556
$ws.Range("F9").Value = 38.3509
$ws.Range("F10").Value = -6.2497
$ws.Range("F11").Value = 52.6723
$ws.Range("F12").Value = -6.932
$ws.Range("F13").Value = 17.5662
$ws.Range("F14").Value = -35.5106
$ws.Range("F15").Value = 0.6286
$ws.Range("F16").Value = -3.1514
$ws.Range("F17").Value = -20.6354
$ws.Range("F18").Value = -0.0175
$ws.Range("F19").Value = -26.9255
$ws.Range("F20").Value = 44.703
$ws.Range("F21").Value = 10.0506
$ws.Range("F22").Value = 84.6016
$ws.Range("F23").Value = -54.4868
$ws.Range("F24").Value = -12.8122
$ws.Range("F25").Value = -9.182700000000001
$ws.Range("F26").Value = 5.9529
$ws.Range("F27").Value = -33.2998
$ws.Range("F28").Value = -20.4441
$ws.Range("F29").Value = -17.1514
$ws.Range("F30").Value = 24.527
$ws.Range("F31").Value = 57.6193
$ws.Range("F32").Value = -1.527
$ws.Range("F33").Value = -5.2378
$ws.Range("F34").Value = 27.4054
$ws.Range("F35").Value = 6.7961
$ws.Range("F36").Value = -5.6683
$ws.Range("F37").Value = 1.4178
$ws.Range("F38").Value = -22.4272
$ws.Range("F39").Value = 12.3741
$ws.Range("F40").Value = -5.138
$ws.Range("F41").Value = -0.1825
$ws.Range("F42").Value = 23.2483
$ws.Range("F43").Value = 14.456
$ws.Range("F44").Value = -11.1739
$ws.Range("F45").Value = 27.112
$ws.Range("F46").Value = -5.6252
$ws.Range("F47").Value = -36.5148
$ws.Range("F48").Value = -27.8397
$ws.Range("F49").Value = -25.4424
$ws.Range("F50").Value = -49.1173
$ws.Range("F51").Value = -51.065
$ws.Range("F52").Value = -35.4517
$ws.Range("F53").Value = -11.9879
$ws.Range("F54").Value = -3.0992
$ws.Range("F55").Value = -15.3441
$ws.Range("F56").Value = -25.937
$ws.Range("F57").Value = -29.1486
$ws.Range("F58").Value = -6.4093
$ws.Range("F59").Value = -23.3046
$ws.Range("F60").Value = -11.2657
$ws.Range("F61").Value = -9.777699999999999
$ws.Range("F62").Value = -16.0561
$ws.Range("F63").Value = -9.932499999999999
$ws.Range("F64").Value = 51.8767
$ws.Range("F65").Value = -43.5191
$ws.Range("F66").Value = 13.7315
$ws.Range("F67").Value = 12.6111
$ws.Range("F68").Value = 31.7532
$ws.Range("F69").Value = -19.9577
$ws.Range("F70").Value = -12.9642
$ws.Range("F71").Value = 13.2432
$ws.Range("F72").Value = 2.8232
$ws.Range("F73").Value = -9.179
$ws.Range("F74").Value = -14.2931
$ws.Range("F75").Value = 28.3699
$ws.Range("F76").Value = 45.5868
